# Remove nchan ("numchans") from the Excel file.
#
# The "numchans" item occupied row 2 of the "config" sheet (with its
# accompanying comment in column D). Deleting that row shifts every
# subsequent row (time, fs, blocktime, filetype, file Precision) up by
# one, and the now-unused "numchans" / its comment strings fall out of
# the shared-string table. The "chanconfig" sheet's own content is
# untouched - it merely ends up pointing at lower shared-string indices.

$wb = $excel.ActiveWorkbook

$config = $wb.Worksheets.Item("config")
$chanconfig = $wb.Worksheets.Item("chanconfig")

# Drop the whole "numchans" row; cells below shift up automatically.
$config.Rows(2).Delete()

# Restore the selection on "config" to where the author last left it,
# without disturbing which sheet is active.
$config.Range("K7").Select()

# Keep "chanconfig" as the active (visible) sheet, as in the source file.
$chanconfig.Activate()
